$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 5: add "Instructor titular:" label and instructor name ---
$ws.Range("D5").Value = "Instructor titular:"
$ws.Range("D5").Font.Bold = $true
$ws.Range("E5").Value = "BERNARDO ZAPATA BAENA"

# --- Schedule grid (rows 15-20): replace "nan" placeholder with actual course name ---
$config = "CONFIGURACION DE SERVICIOS `n BERNARDO ZAPATA BAENA `n 504"
$ingles = "INGLES `n NUEVO INGLÉS `n 504"

$ws.Range("B15").Value = $config
$ws.Range("C15").Value = $config
$ws.Range("D15").Value = $config
$ws.Range("E15").Value = $config
$ws.Range("F15").Value = $ingles

$ws.Range("B16").Value = $config
$ws.Range("C16").Value = $config
$ws.Range("D16").Value = $config
$ws.Range("E16").Value = $config
$ws.Range("F16").Value = $ingles

$ws.Range("B17").Value = $config
$ws.Range("C17").Value = $config
$ws.Range("D17").Value = $config
$ws.Range("E17").Value = $config
$ws.Range("F17").Value = $config

$ws.Range("B18").Value = $config
$ws.Range("C18").Value = $config
$ws.Range("D18").Value = $config
$ws.Range("E18").Value = $config
$ws.Range("F18").Value = $config

$ws.Range("B19").Value = $config
$ws.Range("C19").Value = $config
$ws.Range("D19").Value = $config
$ws.Range("E19").Value = $config
$ws.Range("F19").Value = $config

$ws.Range("B20").Value = $config
$ws.Range("C20").Value = $config
$ws.Range("D20").Value = $config
$ws.Range("E20").Value = $config
$ws.Range("F20").Value = $config
